# Automatische test-sync: 2025-08-03 14:06:50
# Append a new log row (row 6) to the "Logs" sheet, extend the conditional
# formatting ranges that cover the data rows, and bump the Dashboard's
# "Intern verzoek / Actie voor medewerker" count from 4 to 5.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- 1. Append the new row of data (row 6) -------------------------------
$logs.Range("A6").Value2 = "Kun jij dit even regelen?"
$logs.Range("B6").Value2 = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value2 = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D6").Value2 = "Intern verzoek / Actie voor medewerker"

$antwoord = @"
Beste klant,
Bedankt voor uw bericht. Om u beter van dienst te kunnen zijn, hebben wij meer details nodig over wat precies geregeld moet worden. Kunt u specifiek aangeven waarover u hulp nodig heeft?
Met vriendelijke groet,
[Naam bedrijf] E-mailassistent
"@

$logs.Range("E6").Value2 = $antwoord
$logs.Range("F6").Value2 = "2025-08-03 14:06:05"
$logs.Range("G6").Value2 = "Ja"
$logs.Range("H6").Value2 = "Nee"
$logs.Range("I6").Value2 = "Ja"
$logs.Range("J6").Value2 = "Nee"

# The multi-line text in E6 would otherwise trigger an explicit custom row
# height; auto-fit the row back so it serializes like the other data rows
# (no ht/customHeight attribute), matching rows 2-5.
$logs.Rows.Item(6).AutoFit()

# --- 2. Extend the conditional formatting sqref ranges to include row 6 --
$ranges = @("D2:D5", "G2:G5", "H2:H5", "I2:I5", "J2:J5")
foreach ($addr in $ranges) {
    $col = $addr.Substring(0, 1)
    $newRange = $logs.Range($col + "2:" + $col + "6")
    $fcs = $logs.Range($addr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the Dashboard summary count --------------------------------
$dashboard.Range("B2").Value2 = 5
